$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert two new table rows right after the current last data row
#    (row 19). This pushes the old footer rows (24,25) down to (26,27)
#    and creates two fresh blank rows at 20 and 21.
# ------------------------------------------------------------------
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# ------------------------------------------------------------------
# 2. Fix up borders/number formats for the (now) six data rows.
#    Row 19 still carries the old "last row" (closing) border after the
#    insert, so grab that formatting for the new final row (21) first,
#    then stamp the plain "middle" formatting (taken from row 18) onto
#    rows 19 and 20.
# ------------------------------------------------------------------
$ws.Range("B19:J19").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Update the summary figures.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 341640      # VALOR MORA total
$ws.Range("C13").Value = 3           # Cant. Trabajadores
$ws.Range("F13").Value = 2           # Cant. Periodos

# ------------------------------------------------------------------
# 4. Row 16 previously held LUZ DIVINA MADARRIAGA TORRES (22565160);
#    that worker is removed from the data set and replaced with
#    RICARDO's period-2507 record (matching rows 17/18's values).
# ------------------------------------------------------------------
$ws.Range("C16").Value = "73087774"
$ws.Range("D16").Value = "RICARDO ENRIQUE SANJUAN ARANGO"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Rows 17 and 18 (YINA / DORA, period 2507) are unchanged.

# ------------------------------------------------------------------
# 5. New rows 19-21: the same three workers again, now for period 2508.
# ------------------------------------------------------------------
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73087774"
$ws.Range("D19").Value = "RICARDO ENRIQUE SANJUAN ARANGO"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1007280080"
$ws.Range("D20").Value = "YINA PAOLA FERIA MARTINEZ"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047421124"
$ws.Range("D21").Value = "DORA MARTINEZ HERNANDEZ"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500
